$d = $word.ActiveDocument

# Append " – sakavema10@hotmail.com" after "Saira Katherinne Vega Martin"
$d.Content.Find.Execute("Saira Katherinne Vega Martin", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Saira Katherinne Vega Martin – sakavema10@hotmail.com", 2)
